$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.15
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2.38
$ws.Range("I2").Value = 2.68
$ws.Range("Q2").Value = 2.28
$ws.Range("T2").Value = 1.92
$ws.Range("V2").Value = 1.59
$ws.Range("W2").Value = 1.39
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 8.800000000000001
$ws.Range("Z2").Value = 19
$ws.Range("AB2").Value = 12.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AD2").Value = 13.5
$ws.Range("AH2").Value = 25
$ws.Range("AK2").Value = 50

# Row 3
$ws.Range("AF3").Value = 9.800000000000001

# Row 4
$ws.Range("F4").Value = 1.71
$ws.Range("G4").Value = 1.81
$ws.Range("H4").Value = 4.9
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 4.4
$ws.Range("L4").Value = 1.37
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 2.84
$ws.Range("P4").Value = 1.64
$ws.Range("Q4").Value = 1.92
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 3.4
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.17
$ws.Range("W4").Value = 2.22
$ws.Range("AD4").Value = 29
$ws.Range("AI4").Value = 120
